$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","F","G","H","K","M","N","O")
$data = @(
    (0.1973374482087138,0.06454707696524764,0.02648139415447304,0.5088516510678858,0.3506904785613827,0.5205324219275766,0.1857334543312277,0.8891921518417263,1.254872198835618,1.679148079589524),
    (0.1727343585085066,0.0635243486311623,0.02367719059892437,0.5072837530677674,0.3505277101394881,0.5234837231719567,0.1628532198281789,0.7915956472265862,1.270232382984935,1.684763970045054),
    (0.1575943198039909,0.06289570793554589,0.02194330296029534,0.5066281400920474,0.3506671039479983,0.5255064205809461,0.1487305773582079,0.7321862908476504,1.280124697804544,1.689152838062157),
    (0.1514166777537156,0.0626394078455732,0.02123373392115013,0.506438206407104,0.3507840259144857,0.5263836815138561,0.1429573486809943,0.7081002974538961,1.284271791199746,1.691177825620514),
    (0.1503904195227079,0.0625968433204136,0.02111573092899732,0.506411333269412,0.3508070698037642,0.5265325520397894,0.1419976274909658,0.7041081470602819,1.284967411162619,1.691528354982424),
    (0.157511037495496,0.06289225182339919,0.02193374553371541,0.5066252658437804,0.3506684374708655,0.5255180369985695,0.1486527903633998,0.7318609649862111,1.280180157783374,1.689179190301573),
    (0.1888615816316701,0.06419460977469527,0.02551704183539982,0.5082473120274997,0.3505846420168481,0.5215063480093747,0.1778600146948435,0.8554303969021362,1.260072658886677,1.680889161786112),
    (0.2500535586109436,0.0667411532432638,0.03244619269332816,0.513864974203365,0.3523228183286662,0.5153086769171153,0.2345283494013159,1.102113959300965,1.224302244889725,1.672100535724184),
    (0.2948148140778528,0.06860511841703243,0.03747562681186878,0.5194796075458186,0.3547653266031574,0.5117708506320753,0.2757705938645358,1.2864171275899,1.200254585905014,1.670204420159081),
    (0.3151307863412853,0.0694510507959194,0.03974993778005853,0.5223572868608599,0.3561308337237676,0.5103815317605722,0.2944431899662732,1.371017474301098,1.189799852468871,1.670334128177615),
    (0.3228168129561766,0.06977104953577395,0.0406091616600861,0.5234935278167896,0.3566845866557742,0.5098870429047508,0.3015007997066732,1.403169930308678,1.185910677776325,1.670526054501352),
    (0.3211618182771758,0.06970214789696172,0.04042420257691504,0.5232467487365184,0.3565636940437003,0.5099921342318368,0.2999814171850232,1.396240047108449,1.186745174680802,1.670478366304195),
    (0.3157632676799267,0.06947738433796502,0.03982066719476052,0.5224498336356618,0.3561756560130789,0.5103402163893236,0.2950240933628834,1.373660308638023,1.189478488061764,1.670347055264358),
    (0.3124555465726075,0.06933966471506636,0.039450721253246,0.5219677592027452,0.3559427488487472,0.510557543287959,0.2919858416359773,1.359844898480262,1.191161815587936,1.670285225039947),
    (0.2934861414533145,0.06854979027735908,0.03732671704754864,0.5192980551364954,0.3546812137186635,0.5118660715970833,0.2745484645203078,1.280904170482927,1.200947589512479,1.670215920723223),
    (0.2818368114205612,0.06806468095152951,0.03602018798702744,0.5177431478743841,0.3539725190816512,0.512725153089022,0.26382809335405,1.232676055263752,1.207075057954444,1.670427633442813),
    (0.2751321144955909,0.06778547417206227,0.03526742969835084,0.5168792608579551,0.3535888368357405,0.5132399883975154,0.2576537021620737,1.205007624700002,1.210645028058965,1.670642793852451),
    (0.2728612959176075,0.06769090936872857,0.03501234087748628,0.5165919941615371,0.35346303784128,0.5134178613466673,0.2555617473844336,1.195651575515257,1.211861589206321,1.670731679363499),
    (0.2830773526214045,0.06811634115298659,0.03615940281229513,0.5179055184845325,0.3540454825203128,0.5126315587435002,0.2649701598869001,1.237802612079918,1.206418056981743,1.670395429988261),
    (0.317349152445928,0.06954341240778916,0.03999799508834201,0.5226826442128782,0.3562886364485962,0.5102371183904211,0.2964805447309402,1.380289317583149,1.188673752246201,1.670381747725571),
    (0.3397055812810947,0.07047409529773319,0.04249500665254402,0.5260759624230715,0.3579684176979328,0.5088564831499411,0.3169966454649682,1.474092881939654,1.17748375994215,1.671205262082026),
    (0.3277775768885931,0.06997757120442571,0.04116339519293888,0.5242400708457495,0.3570523011726152,0.5095765025033074,0.3060541100483931,1.423963597032227,1.183418797022338,1.670689526158355),
    (0.2825165267706211,0.06809298654098228,0.03609646882358675,0.5178320171329318,0.3540124217623202,0.5126738075279178,0.2644538661622278,1.235484713873916,1.206714940090945,1.670409698122171),
    (0.2335324290917526,0.06605329236406021,0.03058232325219734,0.5120841840942276,0.3516483362561544,0.5168068170535989,0.2192653122546915,1.034871049881147,1.224302244889725,1.673677745072965)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowVals[$j]
    }
}
